# Added pattern fill for negative values
# -> New "Negative Value" checklist row (row 32) on the "BVTs" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BVTs")
$ws.Activate()

$ws.Range("A32").Value = "Negative Value"
$ws.Range("B32").Value = "Check whether chart is ploting for negative values"
$ws.Range("C32").Value = "1. Select value column with negative data to 'Primary Measure' input field`n2. Go to formatting pane.`n3. Turn on Pattern fill toggle."
$ws.Range("D32").Value = "Plot should render for negative values also and data labels should also appear for the respective arc.`nArc for negative values are filled in line pattern."

# Match formatting of the rest of column C (explicit wrap, same as B/D's column default style)
$ws.Range("C32").WrapText = $true

# Row sized to fit the new wrapped text block
$ws.Rows.Item(32).RowHeight = 39.75

# Scroll/select so the new row is in view, like the author's saved state
$excel.ActiveWindow.ScrollRow = 26
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D32").Select()
